# Trade #75 closed at 2026-02-17 12:57:37 - unknown UNKNOWN +0.000%
#
# 1. Summary!B6        Total Trades   74 -> 75
# 2. Summary!B9        Win Rate %     44.59 -> 44
# 3. Strategy Status!D4 Trades        74 -> 75   (MarketMaking row)
# 4. Strategy Status!G4 Win Rate %    44.59 -> 44 (MarketMaking row)
# 5. All Trades     : append new trade row 76
# 6. MarketMaking   : append the same new trade row 76

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1 & 2: Summary sheet totals
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 75
$summary.Range("B9").Value = 44

# ---------------------------------------------------------------------------
# 3 & 4: Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 75
$status.Range("G4").Value = 44

# ---------------------------------------------------------------------------
# 5 & 6: append trade #75 (row 76) to both "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------------
function Add-TradeRow76($ws) {
    $ws.Cells.Item(76, 1).Value = 75

    # Date/time-looking text must stay literal text, not become a date serial.
    $ws.Cells.Item(76, 2).NumberFormat = "@"
    $ws.Cells.Item(76, 2).Value = "2026-02-17"
    $ws.Cells.Item(76, 2).Style = "Normal"

    $ws.Cells.Item(76, 3).Value = "12:57:31"
    $ws.Cells.Item(76, 4).Value = "MarketMaking"
    $ws.Cells.Item(76, 5).Value = "DOWN"
    $ws.Cells.Item(76, 6).Value = 0.42
    $ws.Cells.Item(76, 7).Value = 0.417515
    $ws.Cells.Item(76, 8).Value = "CLOSED"
    $ws.Cells.Item(76, 9).Value = -0.5916
    $ws.Cells.Item(76, 10).Value = -0
    $ws.Cells.Item(76, 11).Value = 100.23
    $ws.Cells.Item(76, 12).Value = 0
    $ws.Cells.Item(76, 13).Value = 0
    $ws.Cells.Item(76, 14).Value = 0.6
    $ws.Cells.Item(76, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(76, 16).Value = "early_exit"
    $ws.Cells.Item(76, 17).Value = 0.11
}

Add-TradeRow76 $wb.Worksheets.Item("All Trades")
Add-TradeRow76 $wb.Worksheets.Item("MarketMaking")
